$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift existing content down by 2 rows, leaving row 2 blank as a
#     separator between the new section title (row 1) and the existing
#     "Description:/MFG Part #:/QTY:" header (now row 3). ---
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(1).Insert()

# --- New section title above the Engine Control Module BOM ---
$ws.Range("A1").Value = "Engine Control Module Parts:"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Underline = $true

# --- New ARDUINO NANO line item appended to the existing BOM table ---
$ws.Range("A27").Value = "ARDUINO NANO"
$ws.Range("B27").Value = "ARDUINO NANO"
$ws.Range("C27").Value = 1

# --- New "Suggested Mechanical Parts" section ---
$ws.Range("A29").Value = "Suggested Mechanical Parts:"
$ws.Range("A29").Font.Bold = $true
$ws.Range("A29").Font.Underline = $true

$ws.Range("A31").Value = "Description:"
$ws.Range("A31").Font.Bold = $true
$ws.Range("B31").Value = "Source:"
$ws.Range("B31").Font.Bold = $true

$ws.Range("A32").Value = "Fuel Injector Barb Adapter"
$ws.Range("B32").Value = "https://www.ebay.com/itm/113711584424"

$ws.Range("A33").Value = "440 cc/min Fuel Injectors"
$ws.Range("B33").Value = "https://www.ebay.com/itm/223283676528"

$ws.Range("A34").Value = "MAP Sensor"
$ws.Range("B34").Value = "https://www.amazon.com/gp/product/B07Z37XG1J"

$ws.Range("A35").Value = "Fuel Filter/Pressure Regulator"
$ws.Range("B35").Value = "https://www.amazon.com/gp/product/B07W9H5TF9"

$ws.Range("A36").Value = "Fuel Pump"
$ws.Range("B36").Value = "https://www.amazon.com/gp/product/B07J39HNTC"

$ws.Range("A37").Value = "Fuel Hose"
$ws.Range("B37").Value = "https://www.amazon.com/gp/product/B00PLKV5H6"

$ws.Range("A38").Value = "High-Energy Ignition Coil (Sold as pack of 8, but only need 1)"
$ws.Range("B38").Value = "https://www.amazon.com/gp/product/B00EOXZLG2"

# --- Column widths widened to fit the new, longer text ---
$ws.Columns.Item(1).ColumnWidth = 75.86328125
$ws.Columns.Item(2).ColumnWidth = 44.9296875

# --- View state: scrolled down, selection parked below the new content ---
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("A40").Select()
